$wb = $excel.ActiveWorkbook

# Update both the "展览" sheet and the "全部类型" sheet, which contain the
# same event rows duplicated across sheets.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F5").Value = 1901
    $ws.Range("F6").Value = 132
}
